$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns to match the shorter text ---
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
